# Updates the cryptos list worksheet with refreshed price/volume data
# (plus a handful of re-ranked coin rows), mirroring the scheduled
# GitHub Actions data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = [ordered]@{
    'D2' = '28.213.89'
    'E2' = '  -0.06%  '
    'D3' = '1.911.73'
    'E3' = '  +2.19%  '
    'D4' = '1.000'
    'E4' = '  -0.25%  '
    'D5' = '314.59'
    'E5' = '  +0.76%  '
    'D6' = '1.000'
    'E6' = '  -0.34%  '
    'D7' = '0.5073'
    'E7' = '  +0.69%  '
    'D8' = '0.3924'
    'E8' = '  +0.55%  '
    'D9' = '0.09348'
    'E9' = '  -1.90%  '
    'D10' = '1.143'
    'E10' = '  +0.30%  '
    'D11' = '41.97'
    'E11' = '  +2.84%  '
    'D12' = '6.405'
    'E12' = '  -0.66%  '
    'B13' = 'Solana'
    'C13' = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
    'D13' = '20.99'
    'E13' = '  +0.14%  '
    'B14' = 'WrappedEther'
    'C14' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D14' = '1.914.72'
    'E14' = '  +2.49%  '
    'D15' = '7.331'
    'E15' = '  -0.90%  '
    'D16' = '1.000'
    'E16' = '  -0.26%  '
    'D17' = '0.00001127'
    'E17' = '  -0.04%  '
    'D18' = '92.60'
    'E18' = '  +0.04%  '
    'D19' = '0.06624'
    'E19' = '  +0.28%  '
    'D20' = '18.06'
    'E20' = '  +2.16%  '
    'E21' = '  -0.24%  '
    'D22' = '6.234'
    'E22' = '  +0.80%  '
    'D23' = '28.268.89'
    'E23' = '  -0.04%  '
    'D24' = '11.56'
    'E24' = '  +2.62%  '
    'D25' = '2.327'
    'E25' = '  +1.43%  '
    'B26' = 'LidoDAOToken'
    'C26' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D26' = '2.592'
    'E26' = '  +0.72%  '
    'B27' = 'WrappedliquidstakedEther2.0'
    'C27' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D27' = '2.141.29'
    'E27' = '  +2.68%  '
    'B28' = 'EthereumClassic'
    'C28' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D28' = '21.21'
    'E28' = '  +0.06%  '
    'B29' = 'Monero'
    'C29' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
    'D29' = '158.26'
    'E29' = '  -0.44%  '
    'B30' = 'BitcoinCash'
    'C30' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D30' = '127.37'
    'E30' = '  -0.17%  '
    'B31' = 'ImmutableX'
    'C31' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
    'D31' = '1.107'
    'E31' = '  +4.09%  '
    'B32' = 'Stellar'
    'C32' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
    'D32' = '0.1076'
    'E32' = '  +1.18%  '
    'B33' = 'Filecoin'
    'C33' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
    'D33' = '5.676'
    'E33' = '  +0.93%  '
    'B34' = 'HuobiToken'
    'C34' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
    'D34' = '3.612'
    'E34' = '  -0.46%  '
    'B35' = 'FraxShare'
    'C35' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D35' = '9.696'
    'E35' = '  +2.27%  '
    'B36' = 'Hedera'
    'C36' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D36' = '0.06702'
    'E36' = '  -0.73%  '
    'B37' = 'VeChain'
    'C37' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
    'D37' = '0.02438'
    'E37' = '  +1.16%  '
    'B38' = 'Algorand'
    'C38' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D38' = '0.2214'
    'E38' = '  +1.17%  '
    'B39' = 'ARBITRUM'
    'C39' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
    'D39' = '1.247'
    'E39' = '  +0.70%  '
    'B40' = 'TrustWalletToken'
    'C40' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
    'D40' = '1.282'
    'E40' = '  +8.16%  '
    'B41' = 'TheSandbox'
    'C41' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
    'D41' = '0.6535'
    'E41' = '  +2.81%  '
    'B42' = 'Aptos'
    'C42' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
    'D42' = '11.58'
    'E42' = '  +0.65%  '
    'B43' = 'InternetComputer(DFINITY)'
    'C43' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
    'D43' = '5.034'
    'E43' = '  +0.65%  '
    'B44' = 'Frax'
    'C44' = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
    'D44' = '1.000'
    'E44' = '  -0.22%  '
    'B45' = 'Decentraland'
    'C45' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
    'D45' = '0.6126'
    'E45' = '  +2.36%  '
    'B46' = 'EnergySwap'
    'C46' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'D46' = '13.35'
    'E46' = '  -1.22%  '
    'B47' = 'PancakeSwap'
    'C47' = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
    'D47' = '3.724'
    'E47' = '  +1.70%  '
    'B48' = 'WEMIXTOKEN'
    'C48' = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
    'D48' = '1.295'
    'E48' = '  +1.24%  '
    'B49' = 'NEARProtocol'
    'C49' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D49' = '2.034'
    'E49' = '  +1.67%  '
    'B50' = 'Quant'
    'C50' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    'D50' = '122.48'
    'E50' = '  -0.76%  '
    'B51' = 'EOS'
    'C51' = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
    'D51' = '1.190'
    'E51' = '  -0.62%  '
}

foreach ($addr in $updates.Keys) {
    $cell = $ws.Range($addr)
    # Source data is plain text (coin names, URLs, price/volume strings
    # such as "1.000" or "28.213.89") and must stay text, not be
    # auto-coerced into numbers by Excel's cell-value parsing.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$addr]
    $cell.Style = "Normal"
}
